$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 4 (shifts existing rows 4:25 down to 5:26),
# matching the weekly update that prepended a newer price observation
# (Fecha 44643) ahead of the existing history.
$ws.Rows("4:4").Insert()

$ws.Range("A4").Value = 6
$ws.Range("B4").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C4").Value = "Metropolitana"
$ws.Range("D4").Value = 44643
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100101
$ws.Range("H4").Value = "Berries"
$ws.Range("I4").Value = 100101006
$ws.Range("J4").Value = "Higo"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "`$/bandeja 7 kilos"
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 2143
$ws.Range("T4").Value = 7
